# Adds the next batch of File Transfer test results (rows 26-37)
# "Data - half done (36 of 72 images)"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 26; Name = "Peter_Jamieson_2015_full_bmp_local_to_remote_from_PI"; B = 3.0295999999999998; C = 2.9251 }
    @{ Row = 27; Name = "Peter_Jamieson_2015_full_bmp_remote_to_local_from_PI"; B = 2.9967000000000001; C = 2.8813 }
    @{ Row = 28; Name = "Peter_Jamieson_2015_full_png_local_to_remote_from_PI"; B = 2.9607999999999999; C = 2.8071999999999999 }
    @{ Row = 29; Name = "Peter_Jamieson_2015_full_png_remote_to_local_from_PI"; B = 2.911; C = 2.8113999999999999 }
    @{ Row = 30; Name = "Peter_Jamieson_2015_half_bmp_local_to_remote_from_PI"; B = 2.9196; C = 2.8525 }
    @{ Row = 31; Name = "Peter_Jamieson_2015_half_bmp_remote_to_local_from_PI"; B = 2.8927; C = 2.7452000000000001 }
    @{ Row = 32; Name = "Peter_Jamieson_2015_half_png_local_to_remote_from_PI"; B = 2.8736000000000002; C = 2.8094999999999999 }
    @{ Row = 33; Name = "Peter_Jamieson_2015_half_png_remote_to_local_from_PI"; B = 2.8765999999999998; C = 2.8129 }
    @{ Row = 34; Name = "Redhawk_logo_Double_bmp_local_to_remote_from_PI"; B = 3.1065; C = 3.0320999999999998 }
    @{ Row = 35; Name = "Redhawk_Logo_double_bmp_remote_to_local_from_PI"; B = 3.1074000000000002; C = 3.0314999999999999 }
    @{ Row = 36; Name = "Redhawk_logo_Double_png_local_to_remote_from_PI"; B = 2.9175; C = 2.8414000000000001 }
    @{ Row = 37; Name = "Redhawk_Logo_double_png_remote_to_local_from_PI"; B = 2.8730000000000002; C = 2.7812999999999999 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A: file/test label (bold, matches the style used by the other rows)
    $ws.Cells.Item($row, 1).Value = $r.Name
    $ws.Cells.Item($row, 1).Font.Bold = $true

    # Columns B/C: measured values
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C

    # Column D: difference formula, continuing the shared formula used by D14:D25
    $ws.Cells.Item($row, 4).Formula = "=B$row-C$row"
}

# Update the view to match: scrolled down with D37 (the last new cell) selected
$ws.Range("D37").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1

